$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status column (E) with the new "Em Andamento" value and corresponding
# progress percentages (F) for the rows that moved from "not started" to
# "in progress".
$ws.Range("E2").Value = "Em Andamento"
$ws.Range("F2").Value = 0.4

$ws.Range("E6").Value = "Em Andamento"
$ws.Range("F6").Value = 0.05

$ws.Range("E7").Value = "Em Andamento"
$ws.Range("F7").Value = 0.05

# Move the active cell selection to D12, matching the saved view state.
$ws.Range("D12").Select()
